# Helper: turn a PowerShell jagged array (array-of-arrays, one per row) into
# a true 2D Variant array that Excel's Range.Value setter expects.
function ConvertTo-Excel2DArray {
    param($rows)
    $nRows = $rows.Count
    $nCols = $rows[0].Count
    $arr = New-Object 'object[,]' $nRows, $nCols
    for ($i = 0; $i -lt $nRows; $i++) {
        for ($j = 0; $j -lt $nCols; $j++) {
            $arr[$i, $j] = $rows[$i][$j]
        }
    }
    return $arr
}

$wb = $excel.ActiveWorkbook

# --- "Generate plots and stats": drop in the newly-collected 0.10 sample
# rows (7-11 and 22-26) on the "0.10" worksheet. The existing AVERAGE()
# and ratio formulas further down (rows 12/14 and 27/29) recompute
# automatically once the source ranges they already cover are populated.
$ws = $wb.Worksheets.Item("0.10")

# 0.10 sheet, 'costs' block (B:L), new sample rows 7-11
$costsNewRows = @(
    @(424994.22190781898, 373216.479737752, 553646.17735026404, 663630.38074962201, 649296.51776430197, 848138.26708517701, 1448680.45223067, 917968.25825960399, 1713406.3151701801, 6042153.9761082996, 14621671.684563201),
    @(424748.18781187001, 427575.34951368801, 548830.69392687595, 670062.56436329102, 643933.26588632702, 853208.06684915302, 1385288.1292256, 936228.82282138604, 1694053.3844773199, 6056190.50548869, 14640138.884045901),
    @(422838.45179500902, 415063.38948454498, 555897.92507237894, 670939.34679340699, 642802.44953961996, 864060.09580281004, 1458620.20081292, 930575.12784951401, 1693101.7454067001, 6007260.4570507696, 14607634.595560201),
    @(442335.63341076399, 407592.37785540399, 526028.57193360105, 649241.04780065303, 637162.18521806097, 835091.55859278701, 1427686.5546063799, 917329.26830926305, 1686294.8904053599, 6039298.6351100001, 14655590.1978737),
    @(430898.01976817899, 398549.51780801802, 534618.03308232501, 658302.19824889803, 655822.11513757904, 856247.37213183998, 1439675.08579146, 890773.59113979305, 1711915.13322142, 5999357.2746209903, 14651346.015616899)
)
$costsNew = ConvertTo-Excel2DArray $costsNewRows
$ws.Range("B7:L11").Value = $costsNew

# 0.10 sheet, 'runtime' block (O:Y), new sample rows 7-11
$runtimeNewRows = @(
    @(11015.0239113718, 20791.9131442904, 19569.887382909601, 24433.4779717028, 22222.119554877201, 26028.685137629502, 29526.020504534201, 39379.541855305397, 72845.273867249402, 180939.75542485699, 559409.67557765497),
    @(11579.7726418823, 19159.459115937301, 20122.922617942, 21321.060145273801, 20517.7532173693, 23763.917727395801, 39186.260439455502, 43750.6276667118, 76132.871767506003, 166978.42688113399, 596200.33151470101),
    @(13017.9940089583, 22304.0826264768, 22075.7972449064, 22728.470731526599, 23600.9952649474, 20735.4797888547, 36941.396327689203, 41761.149210855299, 73986.160583794102, 202753.40924412, 505187.202041968),
    @(10763.1358075886, 17572.94925116, 17789.217315614202, 22576.907206326701, 31850.275343284, 23608.458707108999, 36106.131421402097, 54325.599066913099, 58372.318157926202, 184548.2574068, 574256.10653124703),
    @(13441.3131717592, 15634.484896436299, 19075.265921652299, 24309.171071276, 26080.841124057701, 23263.130672275998, 41612.482778728001, 48734.386064112099, 75153.234601020798, 169167.135344818, 516557.05714039499)
)
$runtimeNew = ConvertTo-Excel2DArray $runtimeNewRows
$ws.Range("O7:Y11").Value = $runtimeNew

# 0.10 sheet, 'costs' block (B:L), new sample rows 22-26
$costsNew2Rows = @(
    @(449840.635895618, 436579.48141684599, 572395.26925963105, 655472.62867309805, 687500.90038117499, 874438.33537245402, 1452671.7333620901, 948177.732654867, 1722605.26161585, 6007004.0678088497, 14588930.785886999),
    @(450184.25380215701, 439537.68802810099, 543050.79489411798, 670000.59476766502, 659351.021731916, 884954.01192013803, 1418719.79018991, 936373.16030338802, 1696464.6736341701, 6004270.9108230304, 14652419.9926401),
    @(442699.01936745999, 400612.82326824701, 554151.71130634204, 669233.78620558104, 675150.80915966094, 865224.005457622, 1450403.2818329299, 925821.02136516199, 1733608.3113192599, 5996150.5952707501, 14588678.7631495),
    @(443738.30427965801, 444239.56713787001, 573609.46217464295, 661128.76748528401, 673163.75734188606, 863487.50474986201, 1462662.7711853599, 913443.21992350498, 1714255.55192142, 6050715.1808049399, 14627096.303513501),
    @(457764.667560993, 425392.856830089, 557221.77865784802, 670958.660737328, 679189.88290436997, 880835.87174135097, 1467021.61452887, 941855.83488003805, 1718786.4889275499, 6028133.4769543102, 14647109.186931901)
)
$costsNew2 = ConvertTo-Excel2DArray $costsNew2Rows
$ws.Range("B22:L26").Value = $costsNew2

# 0.10 sheet, 'runtime' block (O:Y), new sample rows 22-26
$runtimeNew2Rows = @(
    @(8829.8604488372803, 11228.5926844924, 12664.5312011241, 13743.485821411001, 14350.6542555987, 15551.9964639097, 21449.238313361999, 26489.108288660598, 46903.2757859677, 129889.292085543, 419716.14590473397),
    @(8735.2316919714194, 11342.577021569001, 12688.421936705699, 13658.153947442701, 14529.3767638504, 15050.7364273071, 21432.472098618699, 26386.057676747401, 46982.530264183799, 129736.17026954801, 416948.37646745099),
    @(8750.6659757345897, 11279.234528541499, 12667.613415047501, 13763.026200234801, 14776.608467102, 15087.968196719799, 21402.839474380002, 26341.350927948901, 46820.3050885349, 130430.350769311, 417936.79221533198),
    @(8742.8540755063295, 11295.3914664685, 12746.406763792, 13710.479307919701, 14864.016957581, 15078.9930112659, 21487.774129956899, 26463.1190355867, 46898.2236310839, 129164.245016872, 421718.94021890999),
    @(8735.7209958136009, 11345.712447538899, 12789.283905178299, 13679.760305210901, 14727.4562679231, 15040.159927680999, 21429.887272417502, 26434.605209156802, 46722.871776670203, 129578.92547734, 422872.26403504598)
)
$runtimeNew2 = ConvertTo-Excel2DArray $runtimeNew2Rows
$ws.Range("O22:Y26").Value = $runtimeNew2


# --- View-state: make "0.10" the active sheet/tab (it was "0.05" before),
# scroll it so column M is at the left edge, and leave the selection on
# X21 - this also drops tabSelected from the previously-active "0.05" tab.
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 13
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("X21").Select()
